# Diary "Duo Chai" update
# - Row 30 (2/25/2020): Participants grows from "Soobin" to "Soobin, Marc"
# - Row 33 (3/2/2020): homework 3 resubmission entry (was a blank styled row)
# - Row 34 (3/3/2020): homework 5 entry (was a blank styled row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: add Marc as a participant alongside Soobin ---
$ws.Range("C30").Value = "Soobin, Marc"

# --- Row 33: previously blank, styled row -> fill in with homework3 resubmission entry ---
# Copy the formatting (number format / fills / borders / fonts) from the fully
# populated row 30 so the new row matches the rest of the table exactly.
$ws.Range("A30:G30").Copy($ws.Range("A33:G33"))
$ws.Range("A33").Value = 43892
$ws.Range("B33").Value = "20:45-21:20"
$ws.Range("C33").Value = "Soobin, Marc"
$ws.Range("D33").Value = "Finish homework 3 resubmission"
$ws.Range("E33").Value = "Resubmit homework 3"
$ws.Range("F33").Value = "After diving deeper into the key developers, we know more about their developing team members, and their certain roles."
$ws.Range("G33").Value = "Average"

# --- Row 34: previously blank, styled row -> fill in with homework5 entry ---
$ws.Range("A30:G30").Copy($ws.Range("A34:G34"))
$ws.Range("A34").Value = 43893
$ws.Range("B34").Value = "14:00-16:30"
$ws.Range("C34").Value = "Soobin, Marc"
$ws.Range("D34").Value = "Finish homework 5"
$ws.Range("E34").Value = "Finished homework 5"
$ws.Range("F34").Value = "Learnt more about our project by finding design patterns. I feel like patterns are cool and can be signatures of experts."
$ws.Range("G34").Value = "Good"
